# Trade #12 closed at 2026-02-16 22:58:29 - base_strategy DOWN +0.000%
# Append a new trade row (row 13) to both the "All Trades" and "base_strategy"
# worksheets, mirroring the existing OPEN-trade rows already present there.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(13, 1).Value2 = 12                 # A13 Trade #
    # Leading apostrophe forces the date-shaped text to stay literal text
    # (matching the logged value) instead of being auto-converted to a
    # date serial number.
    $ws.Cells.Item(13, 2).Value2 = "'2026-02-16"       # B13 Date
    $ws.Cells.Item(13, 3).Value2 = "22:58:29"          # C13 Time
    $ws.Cells.Item(13, 4).Value2 = "base_strategy"     # D13 Strategy
    $ws.Cells.Item(13, 5).Value2 = "DOWN"              # E13 Side
    $ws.Cells.Item(13, 6).Value2 = 0.5                 # F13 Entry Price
    $ws.Cells.Item(13, 7).Value2 = ""                  # G13 Exit Price (blank)
    $ws.Cells.Item(13, 8).Value2 = "OPEN"              # H13 Status
    $ws.Cells.Item(13, 9).Value2 = 0                   # I13 P&L %
    $ws.Cells.Item(13, 10).Value2 = 0                  # J13 P&L $
    $ws.Cells.Item(13, 11).Value2 = 100                # K13 Capital After
    $ws.Cells.Item(13, 12).Value2 = 0                  # L13 Entry Slippage (bps)
    $ws.Cells.Item(13, 13).Value2 = 0                  # M13 Exit Slippage (bps)
    $ws.Cells.Item(13, 14).Value2 = 0.6                # N13 Confidence
    $ws.Cells.Item(13, 15).Value2 = "Normal spread capture: 19600 bps"  # O13 Entry Reason
    $ws.Cells.Item(13, 16).Value2 = ""                 # P13 Exit Reason (blank)
    $ws.Cells.Item(13, 17).Value2 = 0                  # Q13 Duration (min)
}
